$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 17588
$ws.Range("E2").Value = 2086
$ws.Range("F2").Value = 2086
$ws.Range("G2").Value = 1641
$ws.Range("H2").Value = 1298
$ws.Range("I2").Value = 1285
$ws.Range("J2").Value = 13
$ws.Range("K2").Value = 25273
$ws.Range("L2").Value = 15940
$ws.Range("M2").Value = 9334
$ws.Range("N2").Value = 9280
$ws.Range("O2").Value = 54
$ws.Range("P2").Value = 535
$ws.Range("Q2").Value = 3035
$ws.Range("R2").Value = -3212
$ws.Range("S2").Value = 100
$ws.Range("T2").Value = 3192
$ws.Range("U2").Value = -157
$ws.Range("V2").Value = 11703
$ws.Range("W2").Value = 11.86
$ws.Range("X2").Value = 7.38
$ws.Range("Y2").Value = 14.93
$ws.Range("Z2").Value = 5.38
$ws.Range("AA2").Value = 170.78
$ws.Range("AB2").Value = 1599.52
$ws.Range("AC2").Value = 1255
$ws.Range("AD2").Value = 10.16
$ws.Range("AE2").Value = 9151
$ws.Range("AF2").Value = 1.39
$ws.Range("AG2").Value = 75
$ws.Range("AH2").Value = 0.59
$ws.Range("AI2").Value = 5.94
$ws.Range("AJ2").Value = 96415877

# Row 3
$ws.Range("D3").Value = 18375
$ws.Range("E3").Value = 2249
$ws.Range("F3").Value = 2249
$ws.Range("G3").Value = 1648
$ws.Range("H3").Value = 1271
$ws.Range("I3").Value = 1252
$ws.Range("J3").Value = 19
$ws.Range("K3").Value = 25809
$ws.Range("L3").Value = 15468
$ws.Range("M3").Value = 10340
$ws.Range("N3").Value = 10266
$ws.Range("O3").Value = 74
$ws.Range("P3").Value = 535
$ws.Range("Q3").Value = 3208
$ws.Range("R3").Value = -1421
$ws.Range("S3").Value = -1156
$ws.Range("T3").Value = 1616
$ws.Range("U3").Value = 1593
$ws.Range("V3").Value = 10748
$ws.Range("W3").Value = 12.24
$ws.Range("X3").Value = 6.92
$ws.Range("Y3").Value = 12.81
$ws.Range("Z3").Value = 4.98
$ws.Range("AA3").Value = 149.59
$ws.Range("AB3").Value = 1779.78
$ws.Range("AC3").Value = 1217
$ws.Range("AD3").Value = 10.15
$ws.Range("AE3").Value = 10124
$ws.Range("AF3").Value = 1.22
$ws.Range("AG3").Value = 80
$ws.Range("AH3").Value = 0.65
$ws.Range("AI3").Value = 6.51
$ws.Range("AJ3").Value = 96415877

# Row 4
$ws.Range("D4").Value = 18947
$ws.Range("E4").Value = 2480
$ws.Range("F4").Value = 2480
$ws.Range("G4").Value = 2165
$ws.Range("H4").Value = 1760
$ws.Range("I4").Value = 1738
$ws.Range("J4").Value = 22
$ws.Range("K4").Value = 26179
$ws.Range("L4").Value = 14152
$ws.Range("M4").Value = 12028
$ws.Range("N4").Value = 11932
$ws.Range("O4").Value = 96
$ws.Range("P4").Value = 535
$ws.Range("Q4").Value = 3832
$ws.Range("R4").Value = -1264
$ws.Range("S4").Value = -2374
$ws.Range("T4").Value = 1046
$ws.Range("U4").Value = 2786
$ws.Range("V4").Value = 8375
$ws.Range("W4").Value = 13.09
$ws.Range("X4").Value = 9.29
$ws.Range("Y4").Value = 15.66
$ws.Range("Z4").Value = 6.77
$ws.Range("AA4").Value = 117.66
$ws.Range("AB4").Value = 2094.62
$ws.Range("AC4").Value = 1689
$ws.Range("AD4").Value = 7.7
$ws.Range("AE4").Value = 11766
$ws.Range("AF4").Value = 1.1
$ws.Range("AG4").Value = 100
$ws.Range("AH4").Value = 0.77
$ws.Range("AI4").Value = 5.85
$ws.Range("AJ4").Value = 96415877

# Row 5
$ws.Range("D5").Value = 19648
$ws.Range("E5").Value = 1854
$ws.Range("F5").Value = 1854
$ws.Range("G5").Value = 1667
$ws.Range("H5").Value = 1254
$ws.Range("I5").Value = 1250
$ws.Range("J5").Value = 4
$ws.Range("K5").Value = 27940
$ws.Range("L5").Value = 14551
$ws.Range("M5").Value = 13389
$ws.Range("N5").Value = 13355
$ws.Range("O5").Value = 34
$ws.Range("P5").Value = 541
$ws.Range("Q5").Value = 2551
$ws.Range("R5").Value = -2897
$ws.Range("S5").Value = 260
$ws.Range("T5").Value = 2839
$ws.Range("U5").Value = -288
$ws.Range("V5").Value = 8631
$ws.Range("W5").Value = 9.44
$ws.Range("X5").Value = 6.38
$ws.Range("Y5").Value = 9.88
$ws.Range("Z5").Value = 4.63
$ws.Range("AA5").Value = 108.67
$ws.Range("AB5").Value = 2349.7
$ws.Range("AC5").Value = 1208
$ws.Range("AD5").Value = 9.77
$ws.Range("AE5").Value = 13009
$ws.Range("AF5").Value = 0.91
$ws.Range("AG5").Value = 100
$ws.Range("AH5").Value = 0.85
$ws.Range("AI5").Value = 8.24
$ws.Range("AJ5").Value = 97667877

# Row 6
$ws.Range("D6").Value = 19840
$ws.Range("E6").Value = 1824
$ws.Range("F6").Value = 1824
$ws.Range("G6").Value = 1486
$ws.Range("H6").Value = 1035
$ws.Range("I6").Value = 1036
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 32423
$ws.Range("L6").Value = 18141
$ws.Range("M6").Value = 14282
$ws.Range("N6").Value = 14248
$ws.Range("O6").Value = 34
$ws.Range("P6").Value = 541
$ws.Range("Q6").Value = 1570
$ws.Range("R6").Value = -4516
$ws.Range("S6").Value = 3148
$ws.Range("T6").Value = 4844
$ws.Range("U6").Value = -3274
$ws.Range("V6").Value = 11950
$ws.Range("W6").Value = 9.2
$ws.Range("X6").Value = 5.22
$ws.Range("Y6").Value = 7.5
$ws.Range("Z6").Value = 3.43
$ws.Range("AA6").Value = 127.02
$ws.Range("AB6").Value = 2512.29
$ws.Range("AC6").Value = 994
$ws.Range("AD6").Value = 9.22
$ws.Range("AE6").Value = 13878
$ws.Range("AF6").Value = 0.66
$ws.Range("AG6").Value = 100
$ws.Range("AH6").Value = 1.09
$ws.Range("AI6").Value = 9.94
$ws.Range("AJ6").Value = 97667877

# Row 7
$ws.Range("D7").Value = 20897
$ws.Range("E7").Value = 2198
$ws.Range("G7").Value = 1921
$ws.Range("H7").Value = 1267
$ws.Range("I7").Value = 1266
$ws.Range("K7").Value = 34598
$ws.Range("L7").Value = 19081
$ws.Range("M7").Value = 15517
$ws.Range("N7").Value = 15483
$ws.Range("P7").Value = 540
$ws.Range("Q7").Value = 2775
$ws.Range("R7").Value = -2919
$ws.Range("S7").Value = 388
$ws.Range("T7").Value = 2640
$ws.Range("U7").Value = -226
$ws.Range("W7").Value = 10.52
$ws.Range("X7").Value = 6.06
$ws.Range("Y7").Value = 8.51
$ws.Range("Z7").Value = 3.78
$ws.Range("AA7").Value = 122.97
$ws.Range("AC7").Value = 1215
$ws.Range("AD7").Value = 6.67
$ws.Range("AE7").Value = 15081
$ws.Range("AF7").Value = 0.54
$ws.Range("AG7").Value = 102
$ws.Range("AH7").Value = 1.26
$ws.Range("AI7").Value = 7.88

# Row 8
$ws.Range("D8").Value = 22257
$ws.Range("E8").Value = 2245
$ws.Range("G8").Value = 1955
$ws.Range("H8").Value = 1439
$ws.Range("I8").Value = 1436
$ws.Range("K8").Value = 36048
$ws.Range("L8").Value = 19184
$ws.Range("M8").Value = 16865
$ws.Range("N8").Value = 16827
$ws.Range("P8").Value = 540
$ws.Range("Q8").Value = 3095
$ws.Range("R8").Value = -2350
$ws.Range("S8").Value = -170
$ws.Range("T8").Value = 2364
$ws.Range("U8").Value = 920
$ws.Range("W8").Value = 10.09
$ws.Range("X8").Value = 6.46
$ws.Range("Y8").Value = 8.89
$ws.Range("Z8").Value = 4.07
$ws.Range("AA8").Value = 113.75
$ws.Range("AC8").Value = 1379
$ws.Range("AD8").Value = 5.88
$ws.Range("AE8").Value = 16391
$ws.Range("AF8").Value = 0.49
$ws.Range("AG8").Value = 103
$ws.Range("AH8").Value = 1.27
$ws.Range("AI8").Value = 7.03

# Row 9
$ws.Range("D9").Value = 22977
$ws.Range("E9").Value = 2398
$ws.Range("G9").Value = 2102
$ws.Range("H9").Value = 1539
$ws.Range("I9").Value = 1534
$ws.Range("K9").Value = 37181
$ws.Range("L9").Value = 18871
$ws.Range("M9").Value = 18308
$ws.Range("N9").Value = 18270
$ws.Range("P9").Value = 540
$ws.Range("Q9").Value = 3361
$ws.Range("R9").Value = -2365
$ws.Range("S9").Value = -478
$ws.Range("T9").Value = 2397
$ws.Range("U9").Value = 1198
$ws.Range("W9").Value = 10.44
$ws.Range("X9").Value = 6.7
$ws.Range("Y9").Value = 8.74
$ws.Range("Z9").Value = 4.2
$ws.Range("AA9").Value = 103.07
$ws.Range("AC9").Value = 1473
$ws.Range("AD9").Value = 5.51
$ws.Range("AE9").Value = 17796
$ws.Range("AF9").Value = 0.46
$ws.Range("AG9").Value = 105
$ws.Range("AH9").Value = 1.3
$ws.Range("AI9").Value = 6.7
